$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$emails = @(
    "660r1cw6@gmail.com",
    "a9ausm0a@yahoo.com",
    "sfna3ug2@yahoo.com",
    "68zoqph9@yahoo.com",
    "vgt9znho@gmail.com",
    "re7guys3@hotmail.com",
    "egptfjaz@hotmail.com",
    "46qhzob0@example.com",
    "oicmm8g9@gmail.com",
    "04p9ig59@hotmail.com"
)

$row = 2
foreach ($email in $emails) {
    $ws.Range("E$row").Value = $email
    $row = $row + 1
}
